# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 292
    4  = 11098
    5  = 10312
    6  = 592
    7  = 0
    8  = 727
    9  = 102
    10 = 17
    11 = 31
    13 = 9604
    14 = 14
    16 = 2441
    17 = 38
    18 = 13
    19 = 86
    20 = 391
    21 = 10860
    22 = 10787
    27 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
